$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the updated crypto price/volume values from the scrape.
# NumberFormat is forced to text ("@") before assignment, then the style
# is reset to Normal afterwards, so numeric-looking strings (e.g. "230.51")
# are stored as text -- matching the original inlineStr cell contents --
# without leaving any residual cell formatting behind.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "38.152.77"
$c.Style = "Normal"

$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "  +2.97%  "
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.059.87"
$c.Style = "Normal"

$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "  +2.50%  "
$c.Style = "Normal"

$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "  +0.11%  "
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "230.51"
$c.Style = "Normal"

$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "  +2.08%  "
$c.Style = "Normal"

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.617"
$c.Style = "Normal"

$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "  +1.77%  "
$c.Style = "Normal"

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "59.63"
$c.Style = "Normal"

$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "  +8.39%  "
$c.Style = "Normal"

$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "  +0.00%  "
$c.Style = "Normal"

$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "  +3.42%  "
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.0811"
$c.Style = "Normal"

$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "  +4.42%  "
$c.Style = "Normal"

$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "  +2.65%  "
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "14.79"
$c.Style = "Normal"

$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "  +5.65%  "
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "2.364.81"
$c.Style = "Normal"

$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "  +2.40%  "
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "21.32"
$c.Style = "Normal"

$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "  +8.23%  "
$c.Style = "Normal"

$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "  +2.93%  "
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "5.30"
$c.Style = "Normal"

$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "  +2.30%  "
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "2.066.39"
$c.Style = "Normal"

$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "  +2.82%  "
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "38.095.72"
$c.Style = "Normal"

$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "  +3.02%  "
$c.Style = "Normal"

$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "  +1.35%  "
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "69.95"
$c.Style = "Normal"

$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "  +2.44%  "
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.0₃0835"
$c.Style = "Normal"

$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "  +3.14%  "
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "225.36"
$c.Style = "Normal"

$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "  +0.89%  "
$c.Style = "Normal"

$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "  +0.07%  "
$c.Style = "Normal"

$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "  +0.28%  "
$c.Style = "Normal"

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.25"
$c.Style = "Normal"

$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "  +4.48%  "
$c.Style = "Normal"

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "9.34"
$c.Style = "Normal"

$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "  +4.92%  "
$c.Style = "Normal"

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "166.33"
$c.Style = "Normal"

$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "  +1.34%  "
$c.Style = "Normal"

$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = "  +8.19%  "
$c.Style = "Normal"

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "19.05"
$c.Style = "Normal"

$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = "  +2.44%  "
$c.Style = "Normal"

$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = "  +3.11%  "
$c.Style = "Normal"

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.120"
$c.Style = "Normal"

$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = "  +2.50%  "
$c.Style = "Normal"

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "4.56"
$c.Style = "Normal"

$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = "  +2.79%  "
$c.Style = "Normal"

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "4.60"
$c.Style = "Normal"

$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = "  +2.93%  "
$c.Style = "Normal"

$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = "  +10.54%  "
$c.Style = "Normal"

$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = "  +1.89%  "
$c.Style = "Normal"

$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = "  +1.41%  "
$c.Style = "Normal"

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "6.16"
$c.Style = "Normal"

$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = "  +15.94%  "
$c.Style = "Normal"

$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "  +5.74%  "
$c.Style = "Normal"

$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "  -0.06%  "
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.532.93"
$c.Style = "Normal"

$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "  +5.52%  "
$c.Style = "Normal"

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "98.34"
$c.Style = "Normal"

$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "  +4.17%  "
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.0217"
$c.Style = "Normal"

$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "  +2.96%  "
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "16.94"
$c.Style = "Normal"

$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "  +6.72%  "
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.87"
$c.Style = "Normal"

$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "  +4.28%  "
$c.Style = "Normal"

$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "  +1.93%  "
$c.Style = "Normal"

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.13"
$c.Style = "Normal"

$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "  +1.53%  "
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "4.11"
$c.Style = "Normal"

$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "  +6.01%  "
$c.Style = "Normal"

$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "  +3.05%  "
$c.Style = "Normal"

$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "  +0.44%  "
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "2.252.15"
$c.Style = "Normal"

$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "  +2.40%  "
$c.Style = "Normal"

